$d = $word.ActiveDocument

# 1. Merge "Studerende er interesseret i at " + "programmet er " into a single run.
$d.Content.Find.Execute("at programmet er", $true, $false, $false, $false, $false, $true, 1, $false, "at programmet er", 2) | Out-Null

# 2. Remove the stray _GoBack bookmark that currently sits after "Forudsætninger (Preconditions)".
$d.Bookmarks.Item("_GoBack").Delete()

# 3. "Studerende vælger profil til udregningen." -> "Studerende vælger emne" + " til udregningen." (two runs)
$r = $d.Content
$r.Find.Execute("Studerende vælger profil til udregningen.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Text = "Studerende vælger emne"
$after = $r.Duplicate
$after.Start = $r.End
$after.End = $after.Start
$after.InsertAfter(" til udregningen.")

# 4. "Systemet oplyser hvilken profil der er valgt." -> "Systemet oplyser hvilken emne" + _GoBack bookmark + " der er valgt."
$r2 = $d.Content
$r2.Find.Execute("Systemet oplyser hvilken profil der er valgt.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r2.Text = "Systemet oplyser hvilken emne"
$mid = $r2.Duplicate
$mid.Start = $r2.End
$mid.End = $mid.Start
$d.Bookmarks.Add("_GoBack", $mid) | Out-Null
$tail = $d.Content
$tail.Start = $mid.End
$tail.End = $tail.Start
$tail.InsertAfter(" der er valgt.")

# 5. "Studerende angiver input." -> "Studerende angiver " + "Ft" + " areal" + "."
$r3 = $d.Content
$r3.Find.Execute("Studerende angiver input.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r3.Text = "Studerende angiver "
$p2 = $d.Content
$p2.Start = $r3.End
$p2.End = $p2.Start
$p2.InsertAfter("Ft")
$p3 = $d.Content
$p3.Start = $p2.End
$p3.End = $p3.Start
$p3.InsertAfter(" areal")
$p4 = $d.Content
$p4.Start = $p3.End
$p4.End = $p4.Start
$p4.InsertAfter(".")

# 6. Merge "Systemet viser en meddelelse om " + "at der er skrevet negativt tal." into one run.
$d.Content.Find.Execute("Systemet viser en meddelelse om at der er skrevet negativt tal.", $true, $false, $false, $false, $false, $true, 1, $false, "Systemet viser en meddelelse om at der er skrevet negativt tal.", 2) | Out-Null

# 7. Merge " sekunder i " + "80%" + " af tilfældene." into one run.
$d.Content.Find.Execute("sekunder i 80% af tilfældene.", $true, $false, $false, $false, $false, $true, 1, $false, "sekunder i 80% af tilfældene.", 2) | Out-Null

# 8. Rebuild the "Hyppighed (Frequency of Occurrence)" run structure.
$hr = $d.Content
$hr.Find.Execute("Hyppighed (Frequency of Occurrence)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$hr.Text = "Hyppighed"
$hr.LanguageID = 1033
$hTail = $d.Content
$hTail.Start = $hr.End
$hTail.End = $hTail.Start
$hTail.InsertAfter(" (Frequency of Occurrence)")
$hTail.LanguageID = 1033
